$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tgfb2"
$ws.Range("C2").Value = "Tgfbr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.325336333333333
$ws.Range("H2").Value = 3.976009
$ws.Range("I2").Value = 0.02918077208126263
$ws.Range("J2").Value = 0.02918077208126263
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.21795633333333
$ws.Range("N2").Value = 84.653869
$ws.Range("O2").Value = 0.2124568395711989
$ws.Range("P2").Value = 0.212456839571199
$ws.Range("Q2").Value = 37.39828278098011
$ws.Range("R2").Value = 336.584545028821
$ws.Range("S2").Value = 0.006199654612632535
$ws.Range("T2").Value = 0.006199654612632536

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tgfb2"
$ws.Range("C3").Value = "Tgfbr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.325336333333333
$ws.Range("H3").Value = 3.976009
$ws.Range("I3").Value = 0.02918077208126263
$ws.Range("J3").Value = 0.02918077208126263
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.95730733333333
$ws.Range("N3").Value = 191.871922
$ws.Range("O3").Value = 0.4815432848151522
$ws.Range("P3").Value = 0.4815432848151524
$ws.Range("Q3").Value = 84.7649431910331
$ws.Range("R3").Value = 762.8844887192979
$ws.Range("S3").Value = 0.01405180484145349
$ws.Range("T3").Value = 0.0140518048414535

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tgfb2"
$ws.Range("C4").Value = "Tgfbr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.325336333333333
$ws.Range("H4").Value = 3.976009
$ws.Range("I4").Value = 0.02918077208126263
$ws.Range("J4").Value = 0.02918077208126263
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.47719633333334
$ws.Range("N4").Value = 82.431589
$ws.Range("O4").Value = 0.2068795565595709
$ws.Range("P4").Value = 0.2068795565595709
$ws.Range("Q4").Value = 36.41652663870011
$ws.Range("R4").Value = 327.748739748301
$ws.Range("S4").Value = 0.00603690518823752
$ws.Range("T4").Value = 0.006036905188237521

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Tgfb2"
$ws.Range("C5").Value = "Tgfbr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.325336333333333
$ws.Range("H5").Value = 3.976009
$ws.Range("I5").Value = 0.02918077208126263
$ws.Range("J5").Value = 0.02918077208126263
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 13.164899
$ws.Range("N5").Value = 39.494697
$ws.Range("O5").Value = 0.09912031905407785
$ws.Range("P5").Value = 0.09912031905407789
$ws.Range("Q5").Value = 17.44791896936367
$ws.Range("R5").Value = 157.031270724273
$ws.Range("S5").Value = 0.002892407438939079
$ws.Range("T5").Value = 0.00289240743893908

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tgfb2"
$ws.Range("C6").Value = "Tgfbr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 20.45485233333333
$ws.Range("H6").Value = 61.364557
$ws.Range("I6").Value = 0.4503674794711605
$ws.Range("J6").Value = 0.4503674794711605
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.21795633333333
$ws.Range("N6").Value = 84.653869
$ws.Range("O6").Value = 0.2124568395711989
$ws.Range("P6").Value = 0.212456839571199
$ws.Range("Q6").Value = 577.1941299467813
$ws.Range("R6").Value = 5194.747169521032
$ws.Range("S6").Value = 0.09568365133408956
$ws.Range("T6").Value = 0.09568365133408958

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tgfb2"
$ws.Range("C7").Value = "Tgfbr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 20.45485233333333
$ws.Range("H7").Value = 61.364557
$ws.Range("I7").Value = 0.4503674794711605
$ws.Range("J7").Value = 0.4503674794711605
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 63.95730733333333
$ws.Range("N7").Value = 191.871922
$ws.Range("O7").Value = 0.4815432848151522
$ws.Range("P7").Value = 0.4815432848151524
$ws.Range("Q7").Value = 1308.23727714095
$ws.Range("R7").Value = 11774.13549426855
$ws.Range("S7").Value = 0.2168714354384632
$ws.Range("T7").Value = 0.2168714354384633

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Tgfb2"
$ws.Range("C8").Value = "Tgfbr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 20.45485233333333
$ws.Range("H8").Value = 61.364557
$ws.Range("I8").Value = 0.4503674794711605
$ws.Range("J8").Value = 0.4503674794711605
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 27.47719633333334
$ws.Range("N8").Value = 82.431589
$ws.Range("O8").Value = 0.2068795565595709
$ws.Range("P8").Value = 0.2068795565595709
$ws.Range("Q8").Value = 562.0419935323414
$ws.Range("R8").Value = 5058.377941791073
$ws.Range("S8").Value = 0.09317182444184532
$ws.Range("T8").Value = 0.09317182444184534

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Tgfb2"
$ws.Range("C9").Value = "Tgfbr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 20.45485233333333
$ws.Range("H9").Value = 61.364557
$ws.Range("I9").Value = 0.4503674794711605
$ws.Range("J9").Value = 0.4503674794711605
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 13.164899
$ws.Range("N9").Value = 39.494697
$ws.Range("O9").Value = 0.09912031905407785
$ws.Range("P9").Value = 0.09912031905407789
$ws.Range("Q9").Value = 269.2860650282477
$ws.Range("R9").Value = 2423.574585254229
$ws.Range("S9").Value = 0.04464056825676228
$ws.Range("T9").Value = 0.0446405682567623

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tgfb2"
$ws.Range("C10").Value = "Tgfbr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 23.63794933333334
$ws.Range("H10").Value = 70.913848
$ws.Range("I10").Value = 0.5204517484475769
$ws.Range("J10").Value = 0.5204517484475769
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 28.21795633333333
$ws.Range("N10").Value = 84.653869
$ws.Range("O10").Value = 0.2124568395711989
$ws.Range("P10").Value = 0.212456839571199
$ws.Range("Q10").Value = 667.0146220975458
$ws.Range("R10").Value = 6003.131598877912
$ws.Range("S10").Value = 0.1105735336244768
$ws.Range("T10").Value = 0.1105735336244768

$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Tgfb2"
$ws.Range("C11").Value = "Tgfbr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 23.63794933333334
$ws.Range("H11").Value = 70.913848
$ws.Range("I11").Value = 0.5204517484475769
$ws.Range("J11").Value = 0.5204517484475769
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 63.95730733333333
$ws.Range("N11").Value = 191.871922
$ws.Range("O11").Value = 0.4815432848151522
$ws.Range("P11").Value = 0.4815432848151524
$ws.Range("Q11").Value = 1511.819590241762
$ws.Range("R11").Value = 13606.37631217585
$ws.Range("S11").Value = 0.2506200445352355
$ws.Range("T11").Value = 0.2506200445352356

$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Tgfb2"
$ws.Range("C12").Value = "Tgfbr2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 23.63794933333334
$ws.Range("H12").Value = 70.913848
$ws.Range("I12").Value = 0.5204517484475769
$ws.Range("J12").Value = 0.5204517484475769
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 27.47719633333334
$ws.Range("N12").Value = 82.431589
$ws.Range("O12").Value = 0.2068795565595709
$ws.Range("P12").Value = 0.2068795565595709
$ws.Range("Q12").Value = 649.5045747493858
$ws.Range("R12").Value = 5845.541172744473
$ws.Range("S12").Value = 0.1076708269294881
$ws.Range("T12").Value = 0.1076708269294881

$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Tgfb2"
$ws.Range("C13").Value = "Tgfbr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 23.63794933333334
$ws.Range("H13").Value = 70.913848
$ws.Range("I13").Value = 0.5204517484475769
$ws.Range("J13").Value = 0.5204517484475769
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 13.164899
$ws.Range("N13").Value = 39.494697
$ws.Range("O13").Value = 0.09912031905407785
$ws.Range("P13").Value = 0.09912031905407789
$ws.Range("Q13").Value = 311.1912155404507
$ws.Range("R13").Value = 2800.720939864056
$ws.Range("S13").Value = 0.05158734335837649
$ws.Range("T13").Value = 0.05158734335837651
